$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.155687689781189
$ws.Range("B1").Value = 1.943399786949158
$ws.Range("D1").Value = 1.545423746109009
$ws.Range("E1").Value = 0.8990378975868225
